$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (before) values for the columns that get permuted
# across data rows 2-9: D (Fecha), K (Variedad), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# R (Origen), S (Precio $/Kg).
$cols = @("D","K","L","M","N","O","P","R","S")

$before = @{}
for ($r = 2; $r -le 9; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $before[$r] = $rowVals
}

# Mapping of destination row -> source row (data that should end up there)
$mapping = @{
    2 = 7
    3 = 6
    4 = 9
    5 = 8
    6 = 4
    7 = 3
    8 = 2
    9 = 5
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $before[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
